$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.6753301551942219;   C = 0.3127903958511391;  D = 3.900430680208489;  E = 0.496779210170732;  G = 5.385330441424582 }
    3  = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 26.21740644021617;  E = 0.496779210170732;  G = 31.61296591696135 }
    4  = @{ B = 0.003994804209775715; C = 0.04240448674262143; D = 0.1575252929769615; E = 0.496779210170732;  G = 0.7007037941000906 }
    5  = @{ B = 3.230985683306322;    C = 10.29869402782916;   D = 26.21740644021617;  E = 8.660232485948974;  G = 48.40731863730063 }
    6  = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    7  = @{ B = 0.127881588408715;    C = 0.04240448674262143; D = 3.900430680208489;  E = 0.496779210170732;  G = 4.567495965530558 }
    8  = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732;  G = 9.295990156953671 }
    9  = @{ B = 0.3048080303191223;   C = 0.3127903958511391;  D = 26.21740644021617;  E = 0.496779210170732;  G = 27.33178407655716 }
    10 = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    11 = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    12 = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 26.21740644021617;  E = 0.496779210170732;  G = 31.61296591696135 }
    13 = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 337.1190423067083;  E = 8.660232485948974;  G = 350.6780550592317 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
